# download articles with pandoc title blocks
#
# The original document rendered the title as a Heading1 paragraph
# (wrapped in a bookmark) and the byline as a bold "By Dorothy Day" run.
# Pandoc-style title blocks instead use a dedicated "Title" style for the
# headline and an "Authors" style paragraph containing just the author's
# name - with every word/space of each line broken into its own run
# (this is how pandoc's docx writer emits text), and no bookmark wrapping
# the headline.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-RunsXml([string]$text) {
    # Pandoc emits one <w:r> per word, with the separating spaces as their
    # own runs too - so "Look On The" becomes "Look" / " " / "On" / " " / "The".
    $parts = $text -split ' '
    $xml = New-Object System.Text.StringBuilder
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -gt 0) {
            [void]$xml.Append('<w:r><w:t xml:space="preserve"> </w:t></w:r>')
        }
        [void]$xml.Append('<w:r><w:t xml:space="preserve">')
        [void]$xml.Append($parts[$i])
        [void]$xml.Append('</w:t></w:r>')
    }
    return $xml.ToString()
}

$titleText = "Look On The Face Of Thy Christ"
$authorText = "Dorothy Day"

$titleParaXml = '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + (New-RunsXml $titleText) + '</w:p>'
$authorsParaXml = '<w:p ' + $wordNs + '><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + (New-RunsXml $authorText) + '</w:p>'

# Paragraph 1 is the "Look On The Face Of Thy Christ" heading, paragraph 2
# is the "By Dorothy Day" byline, and paragraph 3 is the first paragraph of
# the article body ("The Catholic Worker, December 1937, ..."). Replace the
# heading+byline span in one shot so it becomes the new Title/Authors pair.
$headingPara = $d.Paragraphs(1)
$bodyPara = $d.Paragraphs(3)
$replaceRange = $d.Range($headingPara.Range.Start, $bodyPara.Range.Start)
$replaceRange.InsertXML($titleParaXml + $authorsParaXml)

# The old title was wrapped in a bookmark
# (<w:bookmarkStart .../>...<w:bookmarkEnd .../>) that is dropped in the
# pandoc-style rendering. The start tag sits right at the top of the
# document; clear it out.
$d.Range(0, 0).InsertXML("")
